$wb = $excel.ActiveWorkbook

# Rename "KS4 (GCSE Units)" -> "KS4-11-12"
$ws = $wb.Worksheets.Item("KS4 (GCSE Units)")
$ws.Name = "KS4-11-12"

# Delete the trailing blank "Sheet1"
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Delete()

# B3:B38 label changed from "Twelve" to "Eleven-Twelve"
$ws.Range("B3:B38").Value = "Eleven-Twelve"

# E21:E38 "Working with X" -> "Programming - Working with X"
$ws.Range("E21:E23").Value = "Programming - Working with Sequence"
$ws.Range("E24:E26").Value = "Programming - Working with Selection"
$ws.Range("E27:E29").Value = "Programming - Working with Iteration"
$ws.Range("E30:E32").Value = "Programming - Working with Subroutines"
$ws.Range("E33").Value = "Programming - Working with Strings"
$ws.Range("E34:E35").Value = "Programming - Working with Lists"
$ws.Range("E36").Value = "Programming - Working with Dictionaries "
$ws.Range("E37:E38").Value = "Programming - Working with Data files"

# Restore the frozen-pane view and move the active selection to E38
$ws.Activate()
$ws.Range("E38").Select()
